$wb = $excel.ActiveWorkbook

# --- Update status text + timestamps (Ready for handoff) ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column(s): "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest Handoff Datetime / Latest HO Xliff Generate Date timestamps
$zhcn.Range("H2").Value = "2016-11-03 20:03:25"
$dede.Range("H2").Value = "2016-11-03 20:03:38"
$overview.Range("G2").Value = "2016-11-03 20:03:38"

# --- Narrow the status columns ---
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(3).ColumnWidth = 16.25
